$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "28.186.54"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "1.915.24"
$ws.Range("E3").Value = "  -3.69%  "
$ws.Range("E4").Value = "  -1.00%  "
$ws.Range("D5").Value = "'327.63"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("D7").Value = "'0.4673"
$ws.Range("E7").Value = "  -6.14%  "
$ws.Range("D8").Value = "'0.4004"
$ws.Range("E8").Value = "  -4.73%  "
$ws.Range("D9").Value = "'53.16"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("D10").Value = "'0.08392"
$ws.Range("E10").Value = "  -9.79%  "
$ws.Range("E11").Value = "  -5.04%  "
$ws.Range("D12").Value = "'22.06"
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("D13").Value = "1.919.72"
$ws.Range("E13").Value = "  -5.56%  "
$ws.Range("D14").Value = "'7.409"
$ws.Range("E14").Value = "  -7.06%  "
$ws.Range("D15").Value = "'6.059"
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "'89.57"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D18").Value = "'0.00001060"
$ws.Range("E18").Value = "  -4.54%  "
$ws.Range("D19").Value = "'0.06592"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").Value = "'17.98"
$ws.Range("E20").Value = "  -7.22%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'5.713"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").Value = "28.179.87"
$ws.Range("E23").Value = "  -3.23%  "
$ws.Range("D24").Value = "'11.31"
$ws.Range("E24").Value = "  -5.83%  "
$ws.Range("D25").Value = "'2.276"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").Value = "2.134.59"
$ws.Range("E26").Value = "  -5.79%  "
$ws.Range("D27").Value = "'153.27"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").Value = "'19.98"
$ws.Range("E28").Value = "  -4.17%  "
$ws.Range("D29").Value = "'2.123"
$ws.Range("E29").Value = "  -6.46%  "
$ws.Range("D30").Value = "'5.712"
$ws.Range("E30").Value = "  -9.81%  "
$ws.Range("D31").Value = "'123.14"
$ws.Range("E31").Value = "  -3.60%  "
$ws.Range("D32").Value = "'0.9720"
$ws.Range("E32").Value = "  -7.65%  "
$ws.Range("D33").Value = "'0.09586"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("D34").Value = "'1.440"
$ws.Range("E34").Value = "  -6.45%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.537"
$ws.Range("E35").Value = "  -4.97%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'3.622"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("D37").Value = "'8.787"
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("D38").Value = "'0.02293"
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("D39").Value = "'0.06158"
$ws.Range("E39").Value = "  -3.87%  "
$ws.Range("D40").Value = "'1.242"
$ws.Range("E40").Value = "  -6.06%  "
$ws.Range("D41").Value = "'0.6118"
$ws.Range("E41").Value = "  -5.86%  "
$ws.Range("D42").Value = "'11.00"
$ws.Range("E42").Value = "  -4.44%  "
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").Value = "'0.1898"
$ws.Range("E44").Value = "  -5.20%  "
$ws.Range("D45").Value = "'1.292"
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.80"
$ws.Range("E46").Value = "  -3.93%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5831"
$ws.Range("E47").Value = "  -6.42%  "
$ws.Range("E48").Value = "  -7.82%  "
$ws.Range("D49").Value = "'3.440"
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").Value = "'0.06870"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").Value = "'109.77"
$ws.Range("E51").Value = "  -3.24%  "
